$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix for table 21 (row 5 / group index 3): add two missing production entries ---
$ws.Range("P5").Value = "309.- Publicado en revista especializada: Uso de cartas de control para el análisis de calidad en manufactura de sacos de polipropileno  Colombia, Revista Biotecnología En El Sector Agropecuario Y Agroindustrial ISSN: 1909-9959, 2006 vol:4 fasc: 1 págs: 67 - 76, DOI:10.1142/30  Autores: SILVIO ANDRES MOSQUERA SANCHEZ, JENNY FLOR CABRERA, JHON JAIRO NARVAEZ"
$ws.Range("V5").Value = "5.- Revista de divulgación : Usos de cartas de control para el análisis de calidad en manufactura de sacos de polipropileno  Chile, Actualidad Iberoamericana ISSN: 0717-3636, 2007 vol:10 fasc: págs: 17 - Autores: SILVIO ANDRES MOSQUERA SANCHEZ, JENNY FLOR CABRERA, JHON JAIRO NARVAEZ"

# --- Fix for table 79 (row 7 / group index 5): add missing scientific consultancy entry ---
$ws.Range("BL7").Value = "19. Consultoría científica: ESTUDIO DE CONTROL ESTADÍSTICO DE LA CALIDAD IMPLEMENTANDO CARTAS DE CONTROL P EN EL PROCESO DE BENEFICIO DEL CAFÉ EN EL PERIODO DE COSECHA ENTRE ABRIL Y MAYO DE 2019, EN LA HACIENDA LOS NARANJOS, UNIDAD PRODUCTIVA SUPRACAFÉ COLOMBIA S.A.  Año de inicio: 2019, Mes de inicio: 4, Año de fin: 2019, Mes de fin: 5  Idioma: Español, Ciudad: POPAYÁN, Disponibilidad: Restringido, Duración: 2, Número del contrato: 039, Institución en la cual prestó el servicio: CORPORACION PARQUE TECNOLOGICO DE INNOVACION TECNICAFE"

# --- Fix for table 81 (column CD, "Trabajos dirigidos/turorías"): fill in missing thesis entries for each group ---
$ws.Range("CD2").Value = "63.- Trabajos de grado de pregrado : Aplicación de Cartas de Control por Atributos, para el monitoreo del numero de unidades defectuosas, en la linea de empaques de polietileno de baja densidad, en la empresa Empaques Transparentes S.A.  Desde 2 2016 hasta Junio 2016, Tipo de orientación: Tutor principal Nombre del estudiante:  Gladys Salas Angulo,  Programa académico: Nivel Técnica Profesional en Operación de Procesos Industriales Número de páginas: 55, Valoración: Distincion meritoria, Institución: Institución Universitaria de Barranquilla  Tutor(es)/Cotutor(es): GENETT ISABEL JIMENEZ DELGADO `n 66.- Trabajos de grado de pregrado : Aplicación de Cartas de Control por Variables, para el Monitoreo del peso de productos en la linea de jamones ahumados de la empresa Alimentos Carnicos S.A.S. - ALICAR  Desde 2 2016 hasta Julio 2016, Tipo de orientación: Tutor principal Nombre del estudiante:  Marcela Patricia Santiago Chavez,  Programa académico: Nivel Técnica Profesional en Operación de Procesos Industriales Número de páginas: 58, Valoración: Distincion meritoria, Institución: Institución Universitaria de Barranquilla  Tutor(es)/Cotutor(es): GENETT ISABEL JIMENEZ DELGADO `n 67.- Trabajos de grado de pregrado : Aplicación de Cartas de Control por Variables, para el Monitoreo del peso de productos en la linea de lomos de atun, de la empresa Gralco S.A.  Desde 2 2016 hasta Julio 2016, Tipo de orientación: Tutor principal Nombre del estudiante:  Joshua Trout Varela,  Programa académico: Nivel Técnica Profesional en Operación de Procesos Industriales Número de páginas: 55, Valoración: Distincion meritoria, Institución: Institución Universitaria de Barranquilla  Tutor(es)/Cotutor(es): GENETT ISABEL JIMENEZ DELGADO"

$ws.Range("CD3").Value = "3.- Trabajos de grado de pregrado : REVISIÓN Y RECOMENDACIONES DE LAS FICHAS TÉCNICAS, ATIQUETAS, CARTAS DE CONTROL Y NORMAS BPM  Desde 8 2018 hasta Diciembre 2018, Tipo de orientación: Tutor principal Nombre del estudiante:  Laura Liesel Riobo Caceres,  Programa académico: Ingeniería Agroindustrial Número de páginas: 64, Valoración: Aprobada, Institución: Universidad Nacional de Colombia - Sede Palmira  Tutor(es)/Cotutor(es): SAUL DUSSAN SARRIA"

$ws.Range("CD4").Value = "5.- Trabajos de grado de pregrado : Implementación de control estadístico de la calidad por medio de cartas de control en el área de extrusión de la empresa Plásticos Correa S.A  Desde 1 2021 hasta Septiembre 2021, Tipo de orientación: Tutor principal Nombre del estudiante:  JESUS DAVID VASQUEZ MONTIEL  Programa académico: INGENIERIA INDUSTRIAL Número de páginas: 55, Valoración: Aprobada, Institución: Plásticos Correa S. A.  Tutor(es)/Cotutor(es): EMIL YESID VELEZ LOPEZ"

$ws.Range("CD6").Value = "90.- Trabajos de grado de pregrado : Diseño de Cartas de Control Para Aseguramiento de la Calidad Analítica del Proceso del Laboratorio de la Empresa VEOLIA aguas de Tunja  Desde 7 2018 hasta Marzo 2019, Tipo de orientación: Coturor/asesor Nombre del estudiante:  Juan Sebastián Parra Chaparro,  Programa académico: Ingenieria Electronica Número de páginas: 0, Valoración: Aprobada, Institución: UNIVERSIDAD SANTO TOMAS SECCIONAL TUNJA  Tutor(es)/Cotutor(es): DANIEL ALEJANDRO RODRIGUEZ CARO , PABLO ANDRES ALVAREZ CAMARGO `n 151.- Trabajos de grado de pregrado : DESARROLLO E IMPLEMENTACIÓN DE UN SOFTWARE PARA EL MANEJO DE LAS CARTAS DE CONTROL Y VERIFICACIÓN DE EQUIPOS EN EL LABORATORIO DE AGUAS DE LA EMPRESA PROACTIVA AGUAS DE TUNJA USANDO EL SOFTWARE LABVIEW  Desde 1 2016 hasta Agosto 2016, Tipo de orientación: Coturor/asesor Nombre del estudiante:  DIANA CAROLINA PLAZAS HERNÁNDEZ,  Programa académico: Ing. Electrónica Número de páginas: 203, Valoración: Distincion meritoria, Institución: UNIVERSIDAD SANTO TOMAS SECCIONAL TUNJA  Tutor(es)/Cotutor(es): CAMILO ERNESTO PARDO BEAINY"

$ws.Range("CD8").Value = "78.- Trabajos de grado de pregrado : Diseño de una herramienta lúdica como apoyo en el proceso enseñanza-aprendizaje de las cartas de control para estudiantes de ingeniería industrial  Desde 2 2018 hasta Mayo 2018, Tipo de orientación: Tutor principal Nombre del estudiante:  Susana Estrada Florez,  Programa académico: Ingeniería Industiral Número de páginas: 46, Valoración: Aprobada, Institución: UNIVERSIDAD DE SAN BUENAVENTURA - SEDE MEDELLÍN  Tutor(es)/Cotutor(es): JONATHAN ESPINOSA"

# Restore default row heights for the rows whose cells now contain embedded newlines,
# so wrapping doesn't leave a stale custom row height behind.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(6).AutoFit()

# --- Mark every data row (2-8) as reviewed in column H ---
$ws.Range("H2").Value = "."
$ws.Range("H3").Value = "."
$ws.Range("H4").Value = "."
$ws.Range("H5").Value = "."
$ws.Range("H6").Value = "."
$ws.Range("H7").Value = "."
$ws.Range("H8").Value = "."

# Leave the view focused near the edited columns, matching the author's session.
$ws.Range("CD1").Select()
